$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-11 (A: id, B: name, C: department, D: reason, E: duration, F: date serial, G: salary)
$data = @(
    @(31589, "Pietro Souza", "Operacoes", "Viagem de negocios", 3, 45082, 8979.370000000001),
    @(61665, "Maria Luísa Rodrigues", "Marketing", "Viagem de negocios", 7, 45097, 9493.059999999999),
    @(13186, "Maria Helena Leão", "Vendas", "Doenca", 5, 45079, 2774.09),
    @(7260, "Heitor Moura", "Recursos Humanos", "Consulta medica", 3, 45081, 6013.48),
    @(60168, "Maria Helena Dias", "Financeiro", "Viagem de negocios", 1, 45088, 3533.53),
    @(37126, "Sra. Natália Fogaça", "Operacoes", "Problemas pessoais", 8, 45088, 2675.16),
    @(45884, "Ágatha Castro", "Recursos Humanos", "Doenca", 8, 45078, 6473.18),
    @(29526, "Marina Cirino", "TI", "Problemas pessoais", 4, 45091, 4812.81),
    @(73212, "Davi Lucca Abreu", "Recursos Humanos", "Problemas pessoais", 5, 45094, 2803.32),
    @(96458, "Isabella Correia", "Financeiro", "Consulta medica", 8, 45093, 3407.57)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r++
}
